$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the report title block (rows 1-5): republic/department/title/subtotal/"as of" rows.
# This promotes the former header row (6) to row 1 and the data rows (7-10) to rows 2-5.
$ws.Rows("1:5").Delete()

# Re-point print area / print titles / filter to the new row 1 header layout.
$ws.PageSetup.PrintArea = '$A$1:$I$5'
$ws.PageSetup.PrintTitleRows = '$1:$1'

$ws.AutoFilterMode = $false
$ws.Range("A1:AI5").AutoFilter()

# The hidden _FilterDatabase name only tracks the header row in the source file;
# row-delete doesn't re-point it automatically, so set it explicitly.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='ALS-CLC 2024'!`$A`$1:`$AI`$1"
    }
}
